# Commit message: "set NM range 1.25-2.25"
# The "Pre-and Post-Test Alerts" sheet's Normal Metabolizer (NM) row
# lists an Activity score range in column B. Update it from "1.5-2.25"
# to "1.25-2.25".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-and Post-Test Alerts")

$ws.Range("B5").Value = "1.25-2.25"
